$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set demscalar_ippu_* rows (cols J:AS) to 1
$ws.Range("J4:AS4").Value = 1
$ws.Range("J5:AS5").Value = 1
$ws.Range("J6:AS6").Value = 1
$ws.Range("J7:AS7").Value = 1
$ws.Range("J8:AS8").Value = 1
$ws.Range("J9:AS9").Value = 1
$ws.Range("J10:AS10").Value = 1
$ws.Range("J12:AS12").Value = 1
$ws.Range("J13:AS13").Value = 1
$ws.Range("J14:AS14").Value = 1
$ws.Range("J15:AS15").Value = 1
$ws.Range("J21:AS21").Value = 1
$ws.Range("J22:AS22").Value = 1

# Set prodinit_ippu_* rows (cols J:AS) to new calibrated trajectory values
$values = @(2787682.15124505, 2822942.68896158, 2608239.23298399, 1855477.74024192, 1921026.44739681, 3242888.7726109, 3242889, 3480461.56613973, 3735438.58990419, 4009095.11390513, 4302799.58979336, 4618020.72136221, 4956334.80897373, 5319433.63203414, 5709132.90893579, 6127381.37677069, 6576270.53621962, 7058045.11034743, 7575114.26960487, 8130063.6791689, 8725668.42886572, 9364906.91033395, 10050975.7108226, 10787305.5981015, 11577578.6764198, 12425746.7993008, 13336051.3312513, 14313044.3572033, 15361611.4457508, 16486996.0800101, 17694825.8782744, 18991140.7355813, 20382423.0269187, 21875630.0231035, 23478228.6814331, 25198232.9850844)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(96, 10 + $i).Value = $values[$i]
}

$values = @(200985.875179983, 171013.883486288, 176181.042500502, 188110.565785273, 173759.385957897, 177210.556906891, 177210.6, 184127.88686801, 191315.184997277, 198783.034080981, 206542.385221553, 214604.616989745, 222981.552110541, 231685.474800395, 240729.148781214, 250125.835997506, 259889.316064136, 270033.90647321, 280574.483589729, 291526.504466774, 302906.029512256, 314729.746040421, 327014.992742691, 339779.785113696, 353042.841869794, 366823.612398818, 381142.30528131, 396019.917925053, 411478.267356369, 427540.022213333, 444228.735987819, 461568.881565121, 479585.88711181, 498306.173364446, 517757.192373816, 537967.467761546)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(97, 10 + $i).Value = $values[$i]
}

$values = @(150675.9084637, 192312.082207071, 237878.153233986, 212413.912472056, 213570.226898081, 160635.280528935, 160635.3, 169573.375788012, 179008.784346544, 188969.198285505, 199483.829975275, 210583.517222113, 222300.813710726, 234670.084479251, 247727.606706649, 261511.676108124, 276062.7192506, 291423.41211765, 307638.805271637, 324756.455980114, 342826.567694017, 361902.137286711, 382039.110485703, 403296.545952913, 425736.788494687, 449425.65190958, 474432.612010159, 500831.010384933, 528698.269498011, 558116.119757351, 589170.839217534, 621953.506620109, 656560.268513606, 693092.62123665, 731657.708591205, 772368.636078945)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(98, 10 + $i).Value = $values[$i]
}

$values = @(11561.7766875378, 10706.6315674835, 9815.45981243169, 12318.5924346722, 9208.43588441504, 9124.09880350419, 9124.099, 9488.58951506494, 9867.6407375019, 10261.8343400579, 10671.7752322072, 11098.0924884146, 11541.440313481, 12002.4990464531, 12481.976204637, 12980.6075693182, 13499.1583148544, 14038.4241828732, 14599.2327033765, 15182.4444646262, 15788.9544337592, 16419.6933301585, 17075.6290536878, 17757.7681699815, 18467.1574550693, 19204.8855017053, 19972.0843898673, 20769.931423989, 21599.650939592, 22462.5161820882, 23359.8512606383, 24293.0331800619, 25263.4939539191, 26272.7228020064, 27322.2684356368, 28413.7414342143)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(99, 10 + $i).Value = $values[$i]
}

$values = @(67663.1409163878, 66711.7067807611, 64788.3152640055, 59019.8420882763, 31763.6507073848, 46865.1823382456, 46865.18, 71584.7586092218, 109342.963478227, 167017.168102328, 255112.295785506, 389674.212540096, 595212.361094555, 909163.971848155, 1388712.97327678, 2121205.61511794, 3240059.93188853, 4949066.83605295, 7559509.10248835, 11546859.1885456, 17637382.9718932, 26940423.6266963, 41150460.1528733, 62855744.0023026, 96009729.6459306, 146651166.619665, 224004012.408135, 342157507.038989, 522632422.359609, 798300909.032367, 1219373911.94495, 1862546716.79404, 2844968420.47977, 4345579759.44834, 6637705821.1924, 10138840156.5743)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(100, 10 + $i).Value = $values[$i]
}

$values = @(369189.428968967, 373643.447799493, 398357.885470165, 443378.387368673, 427988.12578654, 489947.653299795, 489947.7, 530810.866816824, 575082.149239252, 623045.794742111, 675009.757926483, 731307.677767996, 792301.019773802, 858381.396802118, 929973.083445373, 1007535.74012074, 1091567.36435832, 1182607.48823598, 1281240.64249022, 1388100.10954471, 1503871.98955308, 1629299.60556241, 1765188.27608109, 1912410.48569367, 2071911.48692124, 2244715.36929431, 2431931.63460541, 2634762.3205582, 2854509.7185512, 3092584.73514869, 3350515.95092446, 3629959.43484133, 3932709.37718053, 4260709.60929065, 4616066.08411947, 5001060.39766121)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(101, 10 + $i).Value = $values[$i]
}

$values = @(168008.761969233, 140831.880747706, 139970.338995684, 25040.5822645366, 16150.443393682, 44275.9404642896, 44275.94, 49846.2526563996, 56117.3608936506, 63177.431117561, 71125.7218631208, 80073.9792844131, 90148.0082097446, 101489.440849683, 114257.728024521, 128632.380906119, 144815.49479109, 163034.590387425, 183545.812560596, 206637.531510785, 232634.397010686, 261901.902702999, 294851.52462772, 331946.50622252, 373708.371128357, 420724.255364778, 473655.16195903, 533245.254081961, 600332.31734426, 675859.538343936, 760888.765061326, 856615.435531413, 964385.384677786, 1085713.76559793, 1222306.35131689, 1376083.51649368)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(102, 10 + $i).Value = $values[$i]
}

$values = @(290010.37017184, 277279.394364313, 288308.49810505, 267935.955217259, 249831.022617216, 244964.943484318, 244964.9, 254041.753465481, 263454.937845447, 273216.915441342, 283340.610328862, 293839.425468398, 304727.260449482, 316018.529892727, 327728.182533627, 339871.721013492, 352465.222403692, 365525.35949042, 379069.422848116, 393115.343730809, 407681.717811654, 422787.829802115, 438453.678983375, 454700.005683785, 471548.31873739, 489020.923959906, 507140.95367983, 525932.397363787, 545420.133376661, 565629.961918559, 586588.639182201, 608323.91277598, 630864.558459568, 654240.418240716, 678482.439883678, 703622.717881572)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(103, 10 + $i).Value = $values[$i]
}

$values = @(123173.665381735, 121397.047708238, 122577.356115681, 122671.220886475, 122465.803216216, 122503.884767388, 122503.9, 127880.54215903, 133493.162771859, 139352.118828781, 145468.221883827, 151852.758005369, 158517.508602343, 165474.772164531, 172737.386957005, 180318.754710634, 188232.865352345, 196494.322820789, 205118.372015046, 214120.926926096, 223518.600002966, 233328.732807745, 243569.428016039, 254259.582821902, 265418.923808905, 277068.043351675, 289228.437615079, 301922.546221187, 315173.793657184, 329006.632500672, 343446.588542103, 358520.307887615, 374255.606129198, 390681.519672912, 407828.359319874, 425727.7661989)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(104, 10 + $i).Value = $values[$i]
}

$values = @(9774.78152374181, 10141.8210262171, 10841.5240678183, 8813.21189029755, 8950.33407241823, 7327.9264758246, 7327.926, 7244.22517468707, 7161.48039453044, 7079.68073941867, 6998.81541397169, 6918.87374611611, 6839.84518567678, 6761.71930298449, 6684.48578749945, 6608.13444645067, 6532.65520349069, 6458.03809736585, 6384.27328060162, 6311.351018203, 6239.26168636979, 6167.99577122647, 6097.54386756665, 6027.89667761179, 5959.0450097842, 5890.97977749394, 5823.69199793966, 5757.17279092313, 5691.41337767722, 5626.40507970742, 5562.13931764645, 5498.60761012205, 5435.80157263761, 5373.71291646571, 5312.33344755419, 5251.65506544474)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(111, 10 + $i).Value = $values[$i]
}

$values = @(27094.1212550234, 29374.6164374737, 35171.9193004231, 33496.2354516403, 28564.7349286167, 28896.8499895872, 28896.85, 29942.1312384723, 31025.2232718065, 32147.4938239751, 33310.3600934165, 34515.2905426664, 35763.8067527263, 37057.4853445097, 38397.9599697917, 39786.9233741778, 41226.1295346958, 42717.3958747103, 44262.6055589573, 45863.7098715962, 47522.7306802836, 49241.7629893787, 51022.9775855059, 52868.623778815, 54781.0322434002, 56762.6179604654, 58815.8832679521, 60943.4210204812, 63147.9178635983, 65432.1576264587, 67799.0248372343, 70251.508365684, 72792.7051974852, 75425.824345094, 78154.1909000727, 80981.250232001)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(112, 10 + $i).Value = $values[$i]
}

Write-Host "done"